$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -------------------------------------------------
# Row 6 ("urban" row label) in all three languages: Kyrgyz, Russian, English.
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

# Row 7 ("rural" row label) in all three languages: Kyrgyz, Russian, English.
$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Row 2 (A-column / Kyrgyz subtitle): fix spacing/wording of the translation.
$ws.Range("A2").Value = "(жалпы калктын санына карата пайыз менен)"

# --- View / selection state --------------------------------------------
# Scroll the sheet back so column A is visible (drop the old topLeftCell="C1"
# scroll position) and move the active selection to A8.
$null = $ws.Range("A1").Select()
$null = $ws.Range("A8").Select()
